$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 24 (shifts existing rows 24-59 down to 25-60)
$ws.Rows.Item(24).Insert()

# Populate the newly-inserted row 24 with its data
$ws.Cells.Item(24, 1).Value = 11
$ws.Cells.Item(24, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(24, 3).Value = "Bíobío"
$ws.Cells.Item(24, 4).Value = 44987
$ws.Cells.Item(24, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(24, 5).Value = 8
$ws.Cells.Item(24, 6).Value = 100112030
$ws.Cells.Item(24, 7).Value = "Poroto granado"
$ws.Cells.Item(24, 8).Value = "Sin especificar"
$ws.Cells.Item(24, 9).Value = "Primera"
$ws.Cells.Item(24, 10).Value = 35
$ws.Cells.Item(24, 11).Value = 22000
$ws.Cells.Item(24, 12).Value = 23000
$ws.Cells.Item(24, 13).Value = 22429
$ws.Cells.Item(24, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(24, 15).Value = "Región del Maule"
$ws.Cells.Item(24, 16).Value = 897
$ws.Cells.Item(24, 17).Value = 25
$ws.Cells.Item(24, 18).Value = "Hortaliza"
